$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Authors" column (E) for rows whose author-list strings had
# inconsistent/incorrect whitespace cleaned up (per the commit message).
$ws.Range("E3").Value = '[Steven H.%Yale%xref no email%1,          Kejian%Liu%xref no email%1]'
$ws.Range("E4").Value = '[Barret%BP%coreGivesNoEmail%1,                Bartels%CL%coreGivesNoEmail%1,                Bauer%R%coreGivesNoEmail%1,                Brevoort%P%coreGivesNoEmail%1,                Brinkeborn%RM%coreGivesNoEmail%1,                Carr%RJ%coreGivesNoEmail%1,                Ernst%E%coreGivesNoEmail%1,                Giles%JT%coreGivesNoEmail%1,                Hoheisel%D%coreGivesNoEmail%1,                Management%of Influenza in the Southern Hemisphere Trialists Study Group%coreGivesNoEmail%1,                Muller-Jakic%B%coreGivesNoEmail%1,                Stimpel%M%coreGivesNoEmail%1,                Turner%RB%coreGivesNoEmail%1,                Winther%B%coreGivesNoEmail%1]'
$ws.Range("E7").Value = '[G. Frank%Lindenmuth%xref no email%1,          Elise B.%Lindenmuth%xref no email%1]'
$ws.Range("E8").Value = '[ Wolfram%Grimm%null%1,                  Hans-Helge%Müller%null%1,                Wolfram%Grimm%null%1,                Hans-Helge%Müller%null%1]'
$ws.Range("E10").Value = '[H.%Hall%xref no email%1,          M.%Fahlman%xref no email%1,          H.%Engels%xref no email%1]'
$ws.Range("E11").Value = '[Joelle%O’Neil%NULL%1,                 Susan%Hughes%susan.hughes@fresno.ucsf.edu%1,                 Andrea%Lourie%NULL%1,                 John%Zweifler%NULL%1]'
$ws.Range("E12").Value = '[M.%Jawad%NULL%1,                 R.%Schoop%NULL%1,                 A.%Suter%NULL%1,                 P.%Klein%NULL%1,                 R.%Eccles%NULL%1]'
$ws.Range("E13").Value = '[E.%Tiralongo%NULL%1,                 R. A.%Lea%NULL%1,                 S. S.%Wee%NULL%1,                 M. M.%Hanna%NULL%1,                 L. R.%Griffiths%NULL%1]'
$ws.Range("E14").Value = '[Ronald B.%Turner%xref no email%1,          Rudolf%Bauer%xref no email%1,          Karin%Woelkart%xref no email%1,          Thomas C.%Hulsey%xref no email%1,          J. David%Gangemi%xref no email%1]'
$ws.Range("E15").Value = '[Steven J.%Sperber%ssperber@humed.com%1,                 Leena P.%Shah%NULL%1,                 Richard D.%Gilbert%NULL%1,                 Thomas W.%Ritchey%NULL%1,                 Arnold S.%Monto%NULL%1]'
$ws.Range("E18").Value = '[ M.%Dorn%null%1,                  E.%Knick%null%1,                  G.%Lewith%null%1,                M.%Dorn%null%1,                E.%Knick%null%1,                G.%Lewith%null%1]'
